$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = [double]"-8.749027035024054e-08"
$ws.Range("B3").Value = [double]"-5.629767683796403e-08"
$ws.Range("B4").Value = [double]"-8.322058358150745e-08"
$ws.Range("B5").Value = [double]"0.2370952607233896"
$ws.Range("B6").Value = [double]"-8.380401102396733e-08"
$ws.Range("B7").Value = [double]"-9.184626952886123e-08"
$ws.Range("B8").Value = [double]"-1.116335100021377e-07"
$ws.Range("B9").Value = [double]"-9.562249426979451e-08"
$ws.Range("B10").Value = [double]"-6.779877257632571e-08"
$ws.Range("B11").Value = [double]"193.8978870539981"
$ws.Range("B12").Value = [double]"-8.997534466468522e-08"
$ws.Range("B13").Value = [double]"-9.237501327312164e-08"
$ws.Range("B14").Value = [double]"210.203498008616"
$ws.Range("B15").Value = [double]"-8.774419646991013e-08"
$ws.Range("B16").Value = [double]"-8.581867987234468e-08"
$ws.Range("B17").Value = [double]"-8.567352157904525e-08"
$ws.Range("B18").Value = [double]"-8.418472809470153e-08"
$ws.Range("B19").Value = [double]"-8.389669355163533e-08"
$ws.Range("B20").Value = [double]"-8.389667934078061e-08"
$ws.Range("B21").Value = [double]"-414.6210828843942"
$ws.Range("B22").Value = [double]"-8.203294691187802e-08"
$ws.Range("B23").Value = [double]"-8.095615281209005e-08"
$ws.Range("B24").Value = [double]"-7.684110158502897e-08"
$ws.Range("B25").Value = [double]"-7.549891741820107e-08"
$ws.Range("B26").Value = [double]"-7.241886027330097e-08"
$ws.Range("B27").Value = [double]"-7.241886649054991e-08"
$ws.Range("B28").Value = [double]"-7.241886649054991e-08"
$ws.Range("B29").Value = [double]"-7.241886649054991e-08"
$ws.Range("B30").Value = [double]"-7.241886649054991e-08"
$ws.Range("B31").Value = [double]"-7.241886649054991e-08"
$ws.Range("B32").Value = [double]"-7.241886649054991e-08"
$ws.Range("B33").Value = [double]"-7.241886649054991e-08"
$ws.Range("B34").Value = [double]"-7.241886649054991e-08"
$ws.Range("B35").Value = [double]"-7.241886649054991e-08"
$ws.Range("B36").Value = [double]"-7.241886649054991e-08"
$ws.Range("B37").Value = [double]"-7.241886649054991e-08"
$ws.Range("B38").Value = [double]"-7.241886649054991e-08"
$ws.Range("B39").Value = [double]"-7.241886649054991e-08"
$ws.Range("B40").Value = [double]"-7.241886649054991e-08"
$ws.Range("B41").Value = [double]"-7.241886649054991e-08"
$ws.Range("B42").Value = [double]"-7.241886649054991e-08"
$ws.Range("B43").Value = [double]"-7.241886649054991e-08"
$ws.Range("B44").Value = [double]"-7.241886649054991e-08"
$ws.Range("B45").Value = [double]"-7.241886649054991e-08"
$ws.Range("B46").Value = [double]"-7.241886649054991e-08"
$ws.Range("B47").Value = [double]"-7.241886649054991e-08"
$ws.Range("B48").Value = [double]"-7.241886649054991e-08"
$ws.Range("B49").Value = [double]"-7.241886649054991e-08"
$ws.Range("B50").Value = [double]"-7.241886649054991e-08"
$ws.Range("B51").Value = [double]"-7.241886649054991e-08"
$ws.Range("B52").Value = [double]"-7.241886649054991e-08"
$ws.Range("B53").Value = [double]"-7.241886649054991e-08"
$ws.Range("B54").Value = [double]"-7.241886649054991e-08"
$ws.Range("B55").Value = [double]"-7.241886649054991e-08"
$ws.Range("B56").Value = [double]"-7.241886649054991e-08"
$ws.Range("B57").Value = [double]"-7.241886649054991e-08"
$ws.Range("B58").Value = [double]"-7.241886649054991e-08"
$ws.Range("B59").Value = [double]"-7.241886649054991e-08"
$ws.Range("B60").Value = [double]"-7.241886649054991e-08"
$ws.Range("B61").Value = [double]"-7.241886649054991e-08"
$ws.Range("B62").Value = [double]"-7.241886649054991e-08"
